$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed cell values
$ws.Range("B2").Value = 3.4
$ws.Range("C4").Value = 27
$ws.Range("B5").Value = 0.9

# Update the active selection to B5 (was C7)
$ws.Range("B5").Select()
